# Exporting odds to excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Home/Away odds for the COUNTER_STRIKE "KOI vs BLEED" match (row 11)
$ws.Range("E11").Value = 1.75
$ws.Range("F11").Value = 2.0

# Remove the LEAGUE_OF_LEGENDS / EMEA Masters "A One Man Army vs Entropiq" match (row 52).
# Deleting the whole row shifts every following row up by one and prunes the now
# unused "A One Man Army" / "Entropiq" shared strings automatically.
$ws.Rows(52).Delete()

# Update refreshed odds for the LEAGUE_OF_LEGENDS / EMEA Masters "GTZ Esports vs
# Boostgate eSports" match, which is now row 53 after the deletion above.
$ws.Range("E53").Value = 6.0
$ws.Range("F53").Value = 1.1

# Update refreshed odds for the RAINBOW_SIX / EU League "Wolves Esports vs Fnatic"
# match, which is now row 56 after the deletion above.
$ws.Range("E56").Value = 1.92
$ws.Range("F56").Value = 1.8
